# Add a new "Constraints" feature row to the "Features supported" sheet,
# matching the formatting of the existing "Yes" rows, then extend the
# color-scale conditional formatting to cover the new row and refresh the
# active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New row 25: Feature name in column C, "Yes" (Supported) in column E.
$ws.Cells.Item(25, 3).Value = "Constraints"
$ws.Cells.Item(25, 5).Value = "Yes"

# Copy the formatting (bold green font) from an existing "Yes" cell onto the
# new cell so it reuses the same cell style instead of minting a new one.
$ws.Cells.Item(10, 5).Copy()
$ws.Cells.Item(25, 5).PasteSpecial(-4122)

# Extend the 3-color color-scale conditional formatting from E10:E24 to
# E10:E25 so the newly added row is included, recreating the rule so the
# min/max stops go back to "automatic" (no explicit val="0").
$ws.Range("E10:E24").FormatConditions.Delete()
$cs = $ws.Range("E10:E25").FormatConditions.AddColorScale(3)
$cs.ColorScaleCriteria.Item(1).Type = 1
$cs.ColorScaleCriteria.Item(1).FormatColor.Color = 7039736
$cs.ColorScaleCriteria.Item(2).Type = 4
$cs.ColorScaleCriteria.Item(2).Value = 50
$cs.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cs.ColorScaleCriteria.Item(3).Type = 2
$cs.ColorScaleCriteria.Item(3).FormatColor.Color = 8109667

# Move the active selection down to E28, as recorded in the saved view state.
$ws.Range("E28").Select() | Out-Null
